$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume snapshot (coinranking.com scrape).
# Rows 43-51 also shifted: mCoin dropped off the list and every coin
# below it moved up one slot, with EnergySwap newly appended at the end.
$ws.Range("D2").Value = "27.479.60"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "1.617.00"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "211.04"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "22.83"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "0.0611"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.846.92"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "1.616.77"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "65.09"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "27.461.51"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "10.17"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  +5.85%  "
$ws.Range("D25").Value = "150.52"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").Value = "6.84"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "1.457.43"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "0.939"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "0.559"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "'0.860"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'67.60"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'0.990"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.757.12"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.70"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "86.51"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.68"
$ws.Range("E51").Value = "  -0.62%  "
